# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values on the active sheet for rows 2-38 to reflect
# the newly-regenerated strikeout (K) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 9
    3  = 8
    4  = 3
    5  = 7
    6  = 6
    7  = 4
    8  = 6
    9  = 6
    10 = 9
    11 = 4
    12 = 7
    13 = 8
    14 = 5
    15 = 5
    16 = 6
    17 = 5
    18 = 7
    19 = 6
    20 = 9
    21 = 7
    22 = 7
    23 = 9
    24 = 3
    25 = 5
    26 = 7
    27 = 5
    28 = 6
    29 = 6
    30 = 2
    31 = 0
    32 = 3
    33 = 5
    34 = 5
    35 = 5
    36 = 4
    37 = 1
    38 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
